$wb = $excel.ActiveWorkbook

# --- Sheet: Compiled Statement ---
$ws1 = $wb.Worksheets.Item("Compiled Statement")

$row9 = @('restructuring', $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 0, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
for ($c = 1; $c -le 27; $c++) { $ws1.Cells.Item(9, $c).Value = $row9[$c-1] }

$row10 = @('impairment and restructuring', $null, $null, $null, $null, $null, 0, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
for ($c = 1; $c -le 27; $c++) { $ws1.Cells.Item(10, $c).Value = $row10[$c-1] }

$row11 = @('total cost of revenue', 34261, 8278, 11064, 9269, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
for ($c = 1; $c -le 27; $c++) { $ws1.Cells.Item(11, $c).Value = $row11[$c-1] }

$row12 = @('gross margin', 55689, 16260, 17854, 17550, 54457, 72007, 19179, 20048, 20401, 62532, 82933, 22649, 24548, 24046, 72891, 96937, 26152, 28882, 28661, 87195, 115856, 31671, 34768, 33745, 101875, 135620)
for ($c = 1; $c -le 27; $c++) { $ws1.Cells.Item(12, $c).Value = $row12[$c-1] }

$row13 = @('research and development', 13037, 3574, 3504, 3715, 11011, 14726, 3977, 4070, 4316, 12560, 16876, 4565, 4603, 4887, 14382, 19269, 4926, 4899, 5204, 15512, 20716, 5599, 5758, 6306, 18206, 24512)
for ($c = 1; $c -le 27; $c++) { $ws1.Cells.Item(13, $c).Value = $row13[$c-1] }

$row14 = @('sales and marketing', 15539, 3812, 4562, 4335, 13134, 17469, 4098, 4588, 4565, 13648, 18213, 4337, 4933, 4911, 14687, 19598, 4231, 4947, 5082, 15035, 20117, 4547, 5379, 5595, 16230, 21825)
for ($c = 1; $c -le 27; $c++) { $ws1.Cells.Item(14, $c).Value = $row14[$c-1] }

$row15 = @('general and administrative', 4481, 1166, 1109, 1208, 3546, 4754, 1149, 1132, 1179, 3706, 4885, 1061, 1121, 1273, 3838, 5111, 1119, 1139, 1327, 3780, 5107, 1287, 1384, 1480, 4420, 5900)
for ($c = 1; $c -le 27; $c++) { $ws1.Cells.Item(15, $c).Value = $row15[$c-1] }

$row17 = @('other income, net', $null, 276, 490, 349, $null, $null, 266, 127, 145, 584, 729, 0, 194, $null, $null, 77, 248, 440, $null, $null, 1186, 286, 268, $null, $null, 333)
for ($c = 1; $c -le 27; $c++) { $ws1.Cells.Item(17, $c).Value = $row17[$c-1] }

$row18 = @('operating income', 22326, 7708, 8679, 8292, 26766, 35058, 9955, 10258, 10341, 32618, 42959, 12686, 13891, 12975, 39984, 52959, 15876, 17897, 17048, 52868, 69916, 20238, 22247, 20364, 63019, 83383)
for ($c = 1; $c -le 27; $c++) { $ws1.Cells.Item(18, $c).Value = $row18[$c-1] }

$row25 = @('product', $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, ' ')
for ($c = 1; $c -le 27; $c++) { $ws1.Cells.Item(25, $c).Value = $row25[$c-1] }

$row26 = @('diluted', 2.71, 0.84, -0.82, 0.95, 1.18, 2.13, 1.14, 1.08, 1.14, 3.92, 5.06, 1.38, 1.51, 1.4, 4.359999999999999, 5.76, 1.82, 2.03, 2.03, 6.020000000000001, 8.05, 2.71, 2.48, 2.22, 7.43, 9.65)
for ($c = 1; $c -le 27; $c++) { $ws1.Cells.Item(26, $c).Value = $row26[$c-1] }

$row27 = @('revenue', $null, $null, $null, $null, $null, 64497, 17299, 16219, 15448, 50621, 66069, 15768, 18255, 15871, 52170, 68041, 15803, 19460, 16873, 54201, 71074, 16631, 20779, 17366, 55366, 72732)
for ($c = 1; $c -le 27; $c++) { $ws1.Cells.Item(27, $c).Value = $row27[$c-1] }

$row28 = @('weighted average shares outstanding:', $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, ' ')
for ($c = 1; $c -le 27; $c++) { $ws1.Cells.Item(28, $c).Value = $row28[$c-1] }

$row29 = @('cost of revenue', $null, $null, $null, $null, $null, 15420, 3649, 5885, 3441, 12832, 16273, 3305, 4966, 3376, 12641, 16017, 3597, 6058, 4277, 13942, 18219, 3792, 6331, 4584, 14480, 19064)
for ($c = 1; $c -le 27; $c++) { $ws1.Cells.Item(29, $c).Value = $row29[$c-1] }

$row30 = @('basic', 7746, 7708, 7710, 7698, 2, 7700, 7673, 7692, 7672, 1, 7673, 7634, 7621, 7602, 8, 7610, 7566, 7555, 7539, 8, 7547, 7513, 7505, 7493, 3, 7496)
for ($c = 1; $c -le 27; $c++) { $ws1.Cells.Item(30, $c).Value = $row30[$c-1] }

$row31 = @('service and other', $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, ' ')
for ($c = 1; $c -le 27; $c++) { $ws1.Cells.Item(31, $c).Value = $row31[$c-1] }

$row32 = @('revenue', $null, $null, $null, $null, $null, 45863, 11785, 16252, 15123, 44651, 59774, 17287, 18651, 19150, 55824, 74974, 21351, 23616, 24833, 72181, 97014, 28686, 30949, 31994, 93544, 125538)
for ($c = 1; $c -le 27; $c++) { $ws1.Cells.Item(32, $c).Value = $row32[$c-1] }

$row33 = @('diluted', 7832, 7799, 7710, 7794, 0, 7794, 7766, 7768, 7744, 9, 7753, 7710, 7691, 7675, 8, 7683, 7637, 7616, 7597, 11, 7608, 7567, 7555, 7534, 6, 7540)
for ($c = 1; $c -le 27; $c++) { $ws1.Cells.Item(33, $c).Value = $row33[$c-1] }

$row34 = @('cash dividends declared per common share', 1.56, 0.42, 0.42, 0.42, 1.26, 1.68, 0.46, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
for ($c = 1; $c -le 27; $c++) { $ws1.Cells.Item(34, $c).Value = $row34[$c-1] }

# --- Sheet: Separately Compiled Statement ---
$ws2 = $wb.Worksheets.Item("Separately Compiled Statement")

$row2_9 = @('total cost of revenue', 34261, 'impairment and restructuring', $null, $null, $null, $null, 0)
for ($c = 1; $c -le 8; $c++) { $ws2.Cells.Item(9, $c).Value = $row2_9[$c-1] }

$row2_10 = @('gross margin', 55689, 'total cost of revenue', 8278, 11064, 9269, $null, $null)
for ($c = 1; $c -le 8; $c++) { $ws2.Cells.Item(10, $c).Value = $row2_10[$c-1] }

$row2_11 = @('research and development', 13037, 'gross margin', 16260, 17854, 17550, 54457, 72007)
for ($c = 1; $c -le 8; $c++) { $ws2.Cells.Item(11, $c).Value = $row2_11[$c-1] }

$row2_12 = @('sales and marketing', 15539, 'other income , net', $null, $null, $null, $null, 1416)
for ($c = 1; $c -le 8; $c++) { $ws2.Cells.Item(12, $c).Value = $row2_12[$c-1] }

$row2_13 = @('general and administrative', 4481, 'research and development', 3574, 3504, 3715, 11011, 14726)
for ($c = 1; $c -le 8; $c++) { $ws2.Cells.Item(13, $c).Value = $row2_13[$c-1] }

$row2_14 = @('impairment, integration, and restructuring', 306, 'sales and marketing', 3812, 4562, 4335, 13134, 17469)
for ($c = 1; $c -le 8; $c++) { $ws2.Cells.Item(14, $c).Value = $row2_14[$c-1] }

$row2_15 = @('operating income', 22326, 'general and administrative', 1166, 1109, 1208, 3546, 4754)
for ($c = 1; $c -le 8; $c++) { $ws2.Cells.Item(15, $c).Value = $row2_15[$c-1] }

$row2_16 = @('other income , net', 823, 'operating income', 7708, 8679, 8292, 26766, 35058)
for ($c = 1; $c -le 8; $c++) { $ws2.Cells.Item(16, $c).Value = $row2_16[$c-1] }

$row2_24 = @('basic', 7746, 'product', $null, $null, $null, $null, $null)
for ($c = 1; $c -le 8; $c++) { $ws2.Cells.Item(24, $c).Value = $row2_24[$c-1] }

$row2_25 = @('diluted', 7832, 'weighted average shares outstanding:', $null, $null, $null, $null, $null)
for ($c = 1; $c -le 8; $c++) { $ws2.Cells.Item(25, $c).Value = $row2_25[$c-1] }

$row2_26 = @('cash dividends declared per common share', 1.56, 'revenue', $null, $null, $null, $null, 64497)
for ($c = 1; $c -le 8; $c++) { $ws2.Cells.Item(26, $c).Value = $row2_26[$c-1] }

$row2_27 = @($null, $null, 'basic', 7708, 7710, 7698, 2, 7700)
for ($c = 1; $c -le 8; $c++) { $ws2.Cells.Item(27, $c).Value = $row2_27[$c-1] }

$row2_28 = @($null, $null, 'cost of revenue', $null, $null, $null, $null, 15420)
for ($c = 1; $c -le 8; $c++) { $ws2.Cells.Item(28, $c).Value = $row2_28[$c-1] }

$row2_29 = @($null, $null, 'diluted', 7799, 7710, 7794, 0, 7794)
for ($c = 1; $c -le 8; $c++) { $ws2.Cells.Item(29, $c).Value = $row2_29[$c-1] }

$row2_30 = @($null, $null, 'service and other', $null, $null, $null, $null, $null)
for ($c = 1; $c -le 8; $c++) { $ws2.Cells.Item(30, $c).Value = $row2_30[$c-1] }

$row2_31 = @($null, $null, 'cash dividends declared per common share', 0.42, 0.42, 0.42, 1.26, 1.68)
for ($c = 1; $c -le 8; $c++) { $ws2.Cells.Item(31, $c).Value = $row2_31[$c-1] }
